# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, as published in the refreshed gh-pages data dump.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    5  = 8773
    8  = 650
    12 = 9
    17 = 2032
    20 = 314
    22 = 2409
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
